{"js": "// Split the single-run \"Bibliografia\" paragraph into a run that carries\n// explicit line breaks (<w:br/>) before each numbered reference, mirroring\n// the target diff:\n//   ...mentoria.<br/><br/>[1] ...2001.<br/>[2] ...19-27.<br/>[3] ...2004.<br/>[4] ...2019.<br/>[5] ...2019.\n\n// Locate the bibliography paragraph by searching for a stable prefix of its\n// text rather than assuming a fixed paragraph index.\nconst searchResults = context.document.body.search(\n  \"A bibliografia ser\u00e1 recomendada pelos docentes respons\u00e1veis\",\n  { matchCase: false }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the bibliography paragraph.\");\n}\n\nconst targetParagraph = searchResults.items[0].paragraphs.getFirst();\n\n// The five numbered references, each prefixed with its literal \"[n] \" tag,\n// exactly as they appear (split out of) the original single run of text.\nconst intro =\n  \"A bibliografia ser\u00e1 recomendada pelos docentes respons\u00e1veis e obtida na \" +\n  \"busca realizada pelos pr\u00f3prios alunos no in\u00edcio dos projetos. Seguem \" +\n  \"refer\u00eancias no t\u00f3pico de mentoria.\";\n\nconst refs = [\n  \"[1] Peddy, S. The art of mentoring \u2013 Lead, follow and get out of the way. Houston: Bullion Books, 2001.\",\n  \"[2] Zachary, L. J. The Mentor\\u2019s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promo\u00e7\u00e3o do bem-estar em estudantes do ensino superior. In: Programa de Monitoriza\u00e7\u00e3o e Tutorado: oito anos a promover a integra\u00e7\u00e3o e o sucesso acad\u00e9mico no IST. Lisboa: IST Press, 2011. p. 19-27.\",\n  \"[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.\",\n  \"[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.\",\n  \"[5] Diretrizes Curriculares Nacionais para os cursos de gradua\u00e7\u00e3o em Engenharia. Minist\u00e9rio da Educa\u00e7\u00e3o. CNE/CES, 2019.\"\n];\n\n// Word represents a manual line break (<w:br/>) as U+000B inside Range.text.\n// Rebuilding the whole paragraph text in one insertText(\"Replace\") call keeps\n// everything inside a single run \u2014 two breaks separate the intro from the\n// first reference, then one break between each subsequent reference.\nconst newText = intro + \"\\u000b\\u000b\" + refs.join(\"\\u000b\");\n\ntargetParagraph.getRange(\"Whole\").insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# Split the single-run \"Bibliografia\" paragraph into a run that carries\n# explicit line breaks (Chr(11) -> <w:br/>) before each numbered reference,\n# mirroring the target diff:\n#   ...mentoria.<br/><br/>[1] ...2001.<br/>[2] ...19-27.<br/>[3] ...2004.<br/>[4] ...2019.<br/>[5] ...2019.\n\n$d = $word.ActiveDocument\n\n# Locate the bibliography paragraph by searching for a stable prefix of its\n# text rather than assuming a fixed paragraph index.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"A bibliografia ser\u00e1 recomendada pelos docentes respons\u00e1veis\")\nif (-not $found) {\n    throw \"Could not find the bibliography paragraph.\"\n}\n\n# Expand the (small) find hit out to the whole enclosing paragraph.\n$rng.Expand(4)  # wdParagraph\n\n$br = [char]11  # manual line break, serializes as <w:br/>\n\n$intro = \"A bibliografia ser\u00e1 recomendada pelos docentes respons\u00e1veis e obtida na busca realizada pelos pr\u00f3prios alunos no in\u00edcio dos projetos. Seguem refer\u00eancias no t\u00f3pico de mentoria.\"\n\n$ref1 = \"[1] Peddy, S. The art of mentoring \u2013 Lead, follow and get out of the way. Houston: Bullion Books, 2001.\"\n$ref2 = \"[2] Zachary, L. J. The Mentor\u2019s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promo\u00e7\u00e3o do bem-estar em estudantes do ensino superior. In: Programa de Monitoriza\u00e7\u00e3o e Tutorado: oito anos a promover a integra\u00e7\u00e3o e o sucesso acad\u00e9mico no IST. Lisboa: IST Press, 2011. p. 19-27.\"\n$ref3 = \"[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.\"\n$ref4 = \"[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.\"\n$ref5 = \"[5] Diretrizes Curriculares Nacionais para os cursos de gradua\u00e7\u00e3o em Engenharia. Minist\u00e9rio da Educa\u00e7\u00e3o. CNE/CES, 2019.\"\n\n# Word represents a manual line break (<w:br/>) as Chr(11) inside Range.Text.\n# Assigning the whole paragraph text in one shot keeps everything inside a\n# single run \u2014 two breaks separate the intro from the first reference, then\n# one break between each subsequent reference.\n$rng.Text = $intro + $br + $br + $ref1 + $br + $ref2 + $br + $ref3 + $br + $ref4 + $br + $ref5\n"}
